$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = [double]"260.0002878026372"
$ws.Cells.Item(2, 2).Value = [double]"249.9983883126787"
$ws.Cells.Item(2, 3).Value = [double]"10.864937702144982"
$ws.Cells.Item(2, 4).Value = [double]"0.9999992802068782"
$ws.Cells.Item(2, 5).Value = [double]"0.0023121071142810433"
$ws.Cells.Item(2, 6).Value = [double]"0.00041287812123195176"
$ws.Cells.Item(2, 7).Value = [double]"4.3459450592925667e-10"
$ws.Cells.Item(3, 1).Value = [double]"259.9943134574173"
$ws.Cells.Item(3, 2).Value = [double]"249.98096178032264"
$ws.Cells.Item(3, 3).Value = [double]"11.650663043343412"
$ws.Cells.Item(3, 4).Value = [double]"0.9996347288572167"
$ws.Cells.Item(3, 5).Value = [double]"0.025991224768841383"
$ws.Cells.Item(3, 6).Value = [double]"-0.008896521975065502"
$ws.Cells.Item(3, 7).Value = [double]"-2.464209848053743e-5"
$ws.Cells.Item(4, 1).Value = [double]"259.98780662654"
$ws.Cells.Item(4, 2).Value = [double]"249.96186399412224"
$ws.Cells.Item(4, 3).Value = [double]"12.364598566779904"
$ws.Cells.Item(4, 4).Value = [double]"0.9995256232158286"
$ws.Cells.Item(4, 5).Value = [double]"0.029542263969808663"
$ws.Cells.Item(4, 6).Value = [double]"-0.010044562587369039"
$ws.Cells.Item(4, 7).Value = [double]"9.872360942379693e-5"
$ws.Cells.Item(5, 1).Value = [double]"259.9809551625287"
$ws.Cells.Item(5, 2).Value = [double]"249.9417658079292"
$ws.Cells.Item(5, 3).Value = [double]"13.013217337176792"
$ws.Cells.Item(5, 4).Value = [double]"0.9994356122958158"
$ws.Cells.Item(5, 5).Value = [double]"0.03217757377601739"
$ws.Cells.Item(5, 6).Value = [double]"-0.010940528025014277"
$ws.Cells.Item(5, 7).Value = [double]"0.00012159830154739718"
$ws.Cells.Item(6, 1).Value = [double]"259.97190486245125"
$ws.Cells.Item(6, 2).Value = [double]"249.9214671887811"
$ws.Cells.Item(6, 3).Value = [double]"13.601510835147828"
$ws.Cells.Item(6, 4).Value = [double]"0.9993408398482525"
$ws.Cells.Item(6, 5).Value = [double]"0.03358208955349479"
$ws.Cells.Item(6, 6).Value = [double]"-0.014930304133917233"
$ws.Cells.Item(6, 7).Value = [double]"0.0015697737994939363"
$ws.Cells.Item(7, 1).Value = [double]"259.9608918008956"
$ws.Cells.Item(7, 2).Value = [double]"249.8965285064425"
$ws.Cells.Item(7, 3).Value = [double]"14.140249845106236"
$ws.Cells.Item(7, 4).Value = [double]"0.9989252861421347"
$ws.Cells.Item(7, 5).Value = [double]"0.04236554262252619"
$ws.Cells.Item(7, 6).Value = [double]"-0.018709350646830058"
$ws.Cells.Item(7, 7).Value = [double]"0.00017885793907870264"
$ws.Cells.Item(8, 1).Value = [double]"259.94948017549103"
$ws.Cells.Item(8, 2).Value = [double]"249.87069881230076"
$ws.Cells.Item(8, 3).Value = [double]"14.629993586820746"
$ws.Cells.Item(8, 4).Value = [double]"0.9987797543469268"
$ws.Cells.Item(8, 5).Value = [double]"0.045128713646074486"
$ws.Cells.Item(8, 6).Value = [double]"-0.019938564893863475"
$ws.Cells.Item(8, 7).Value = [double]"0.00028054153313486843"
$ws.Cells.Item(9, 1).Value = [double]"259.9376683411087"
$ws.Cells.Item(9, 2).Value = [double]"249.8439818388795"
$ws.Cells.Item(9, 3).Value = [double]"15.075284661569924"
$ws.Cells.Item(9, 4).Value = [double]"0.9986232982922177"
$ws.Cells.Item(9, 5).Value = [double]"0.04791653358545916"
$ws.Cells.Item(9, 6).Value = [double]"-0.021185065185635486"
$ws.Cells.Item(9, 7).Value = [double]"0.00044561234538955015"
$ws.Cells.Item(10, 1).Value = [double]"259.92571193489925"
$ws.Cells.Item(10, 2).Value = [double]"249.8168792518992"
$ws.Cells.Item(10, 3).Value = [double]"15.482969487005084"
$ws.Cells.Item(10, 4).Value = [double]"0.9985227276614603"
$ws.Cells.Item(10, 5).Value = [double]"0.04972491671783179"
$ws.Cells.Item(10, 6).Value = [double]"-0.021936027509565623"
$ws.Cells.Item(10, 7).Value = [double]"2.238909846392911e-6"
$ws.Cells.Item(11, 1).Value = [double]"259.9133446874144"
$ws.Cells.Item(11, 2).Value = [double]"249.78884544578742"
$ws.Cells.Item(11, 3).Value = [double]"15.854173972949564"
$ws.Cells.Item(11, 4).Value = [double]"0.9983464869872377"
$ws.Cells.Item(11, 5).Value = [double]"0.05260333860031909"
$ws.Cells.Item(11, 6).Value = [double]"-0.023205881300754166"
$ws.Cells.Item(11, 7).Value = [double]"2.6295146670099873e-6"
$ws.Cells.Item(12, 1).Value = [double]"259.9005723990497"
$ws.Cells.Item(12, 2).Value = [double]"249.75989357111655"
$ws.Cells.Item(12, 3).Value = [double]"16.192295387349542"
$ws.Cells.Item(12, 4).Value = [double]"0.9981610841378784"
$ws.Cells.Item(12, 5).Value = [double]"0.055469756535581956"
$ws.Cells.Item(12, 6).Value = [double]"-0.024470442748811983"
$ws.Cells.Item(12, 7).Value = [double]"3.0534625645564738e-6"
$ws.Cells.Item(13, 1).Value = [double]"259.88740527085"
$ws.Cells.Item(13, 2).Value = [double]"249.730046750286"
$ws.Cells.Item(13, 3).Value = [double]"16.50040697426785"
$ws.Cells.Item(13, 4).Value = [double]"0.9979685181365461"
$ws.Cells.Item(13, 5).Value = [double]"0.05829723510096724"
$ws.Cells.Item(13, 6).Value = [double]"-0.025717826580947716"
$ws.Cells.Item(13, 7).Value = [double]"3.502866571498843e-6"
$ws.Cells.Item(14, 1).Value = [double]"259.8738589778232"
$ws.Cells.Item(14, 2).Value = [double]"249.69934051227284"
$ws.Cells.Item(14, 3).Value = [double]"16.78128437608075"
$ws.Cells.Item(14, 4).Value = [double]"0.9977715286412769"
$ws.Cells.Item(14, 5).Value = [double]"0.061053623116909386"
$ws.Cells.Item(14, 6).Value = [double]"-0.026933847659892544"
$ws.Cells.Item(14, 7).Value = [double]"3.977959517857644e-6"
$ws.Cells.Item(15, 1).Value = [double]"259.8599538086463"
$ws.Cells.Item(15, 2).Value = [double]"249.66782083975258"
$ws.Cells.Item(15, 3).Value = [double]"17.037430959878854"
$ws.Cells.Item(15, 4).Value = [double]"0.9975730117553784"
$ws.Cells.Item(15, 5).Value = [double]"0.06371019869544446"
$ws.Cells.Item(15, 6).Value = [double]"-0.028105831400714102"
$ws.Cells.Item(15, 7).Value = [double]"4.468677839944895e-6"
$ws.Cells.Item(16, 1).Value = [double]"259.8457163063952"
$ws.Cells.Item(16, 2).Value = [double]"249.63554788764233"
$ws.Cells.Item(16, 3).Value = [double]"17.271099061342262"
$ws.Cells.Item(16, 4).Value = [double]"0.9973767992148367"
$ws.Cells.Item(16, 5).Value = [double]"0.0662306851274869"
$ws.Cells.Item(16, 6).Value = [double]"-0.029217769797641024"
$ws.Cells.Item(16, 7).Value = [double]"4.964730085580194e-6"
$ws.Cells.Item(17, 1).Value = [double]"259.8311794244822"
$ws.Cells.Item(17, 2).Value = [double]"249.60259633542256"
$ws.Cells.Item(17, 3).Value = [double]"17.48431022869064"
$ws.Cells.Item(17, 4).Value = [double]"0.9971871889854145"
$ws.Cells.Item(17, 5).Value = [double]"0.06857783831390374"
$ws.Cells.Item(17, 6).Value = [double]"-0.030253226945523157"
$ws.Cells.Item(17, 7).Value = [double]"5.451735191867608e-6"
$ws.Cells.Item(18, 1).Value = [double]"259.8163828791311"
$ws.Cells.Item(18, 2).Value = [double]"249.56905618570704"
$ws.Cells.Item(18, 3).Value = [double]"17.678873797174923"
$ws.Cells.Item(18, 4).Value = [double]"0.9970089581582066"
$ws.Cells.Item(18, 5).Value = [double]"0.07071264244784971"
$ws.Cells.Item(18, 6).Value = [double]"-0.031194984801632405"
$ws.Cells.Item(18, 7).Value = [double]"5.913288967981186e-6"
$ws.Cells.Item(19, 1).Value = [double]"259.8013735517305"
$ws.Cells.Item(19, 2).Value = [double]"249.53503367670643"
$ws.Cells.Item(19, 3).Value = [double]"17.85640416313118"
$ws.Cells.Item(19, 4).Value = [double]"0.9968473161545472"
$ws.Cells.Item(19, 5).Value = [double]"0.07259408695490939"
$ws.Cells.Item(19, 6).Value = [double]"-0.03202494547770348"
$ws.Cells.Item(19, 7).Value = [double]"6.332414334985943e-6"
$ws.Cells.Item(20, 1).Value = [double]"259.78620569262"
$ws.Cells.Item(20, 2).Value = [double]"249.50065174152684"
$ws.Cells.Item(20, 3).Value = [double]"18.018337157208496"
$ws.Cells.Item(20, 4).Value = [double]"0.9967077347576488"
$ws.Cells.Item(20, 5).Value = [double]"0.07418009394848103"
$ws.Cells.Item(20, 6).Value = [double]"-0.03272453807429451"
$ws.Cells.Item(20, 7).Value = [double]"6.6818393170523605e-6"
$ws.Cells.Item(21, 1).Value = [double]"259.7709412632828"
$ws.Cells.Item(21, 2).Value = [double]"249.46605078714907"
$ws.Cells.Item(21, 3).Value = [double]"18.165945706255197"
$ws.Cells.Item(21, 4).Value = [double]"0.9965958831017835"
$ws.Cells.Item(21, 5).Value = [double]"0.07542678174985988"
$ws.Cells.Item(21, 6).Value = [double]"-0.03327439850876786"
$ws.Cells.Item(21, 7).Value = [double]"6.94007372086593e-6"
$ws.Cells.Item(22, 1).Value = [double]"259.7597397150033"
$ws.Cells.Item(22, 2).Value = [double]"249.35056004508476"
$ws.Cells.Item(22, 3).Value = [double]"18.443748366470636"
$ws.Cells.Item(22, 4).Value = [double]"0.9734502820551176"
$ws.Cells.Item(22, 5).Value = [double]"0.22730960021155436"
$ws.Cells.Item(22, 6).Value = [double]"-0.022049355369581358"
$ws.Cells.Item(22, 7).Value = [double]"0.002889459305219203"
$ws.Cells.Item(23, 1).Value = [double]"259.74796131082417"
$ws.Cells.Item(23, 2).Value = [double]"249.2336738917258"
$ws.Cells.Item(23, 3).Value = [double]"18.705956922120883"
$ws.Cells.Item(23, 4).Value = [double]"0.9722191723450032"
$ws.Cells.Item(23, 5).Value = [double]"0.23232227975604727"
$ws.Cells.Item(23, 6).Value = [double]"-0.023414118905147004"
$ws.Cells.Item(23, 7).Value = [double]"0.004192637158709531"
$ws.Cells.Item(24, 1).Value = [double]"259.7349122729767"
$ws.Cells.Item(24, 2).Value = [double]"249.11533260215643"
$ws.Cells.Item(24, 3).Value = [double]"18.954244548322357"
$ws.Cells.Item(24, 4).Value = [double]"0.9709114025218826"
$ws.Cells.Item(24, 5).Value = [double]"0.23727647614345052"
$ws.Cells.Item(24, 6).Value = [double]"-0.026168808041649876"
$ws.Cells.Item(24, 7).Value = [double]"0.008368521345738499"
$ws.Cells.Item(25, 1).Value = [double]"259.72067894898225"
$ws.Cells.Item(25, 2).Value = [double]"248.9954547656741"
$ws.Cells.Item(25, 3).Value = [double]"19.1907966531084"
$ws.Cells.Item(25, 4).Value = [double]"0.9695879288138133"
$ws.Cells.Item(25, 5).Value = [double]"0.24212118322699067"
$ws.Cells.Item(25, 6).Value = [double]"-0.028754599482459897"
$ws.Cells.Item(25, 7).Value = [double]"0.01216024144421478"
$ws.Cells.Item(26, 1).Value = [double]"259.7046687860051"
$ws.Cells.Item(26, 2).Value = [double]"248.87393983900384"
$ws.Cells.Item(26, 3).Value = [double]"19.417563883320582"
$ws.Cells.Item(26, 4).Value = [double]"0.9681495902146381"
$ws.Cells.Item(26, 5).Value = [double]"0.24689374402560413"
$ws.Cells.Item(26, 6).Value = [double]"-0.03253853126279768"
$ws.Cells.Item(26, 7).Value = [double]"0.018374363390541586"
$ws.Cells.Item(27, 1).Value = [double]"259.6862900810947"
$ws.Cells.Item(27, 2).Value = [double]"248.7538857734578"
$ws.Cells.Item(27, 3).Value = [double]"19.62702771525995"
$ws.Cells.Item(27, 4).Value = [double]"0.9677547938252662"
$ws.Cells.Item(27, 5).Value = [double]"0.24680026071070504"
$ws.Cells.Item(27, 6).Value = [double]"-0.037795349347847396"
$ws.Cells.Item(27, 7).Value = [double]"0.02769372374501736"
$ws.Cells.Item(28, 1).Value = [double]"259.66597569972646"
$ws.Cells.Item(28, 2).Value = [double]"248.63520290076278"
$ws.Cells.Item(28, 3).Value = [double]"19.821257110057036"
$ws.Cells.Item(28, 4).Value = [double]"0.9673485350075756"
$ws.Cells.Item(28, 5).Value = [double]"0.24652716628489404"
$ws.Cells.Item(28, 6).Value = [double]"-0.042214181239455716"
$ws.Cells.Item(28, 7).Value = [double]"0.03635371988019857"
$ws.Cells.Item(29, 1).Value = [double]"259.6415924903072"
$ws.Cells.Item(29, 2).Value = [double]"248.52297934142803"
$ws.Cells.Item(29, 3).Value = [double]"19.94159831372103"
$ws.Cells.Item(29, 4).Value = [double]"0.9659631626059105"
$ws.Cells.Item(29, 5).Value = [double]"0.24565657468854482"
$ws.Cells.Item(29, 6).Value = [double]"-0.05340003176536985"
$ws.Cells.Item(29, 7).Value = [double]"0.05756703853948673"
$ws.Cells.Item(30, 1).Value = [double]"259.62376052655884"
$ws.Cells.Item(30, 2).Value = [double]"248.41794482144354"
$ws.Cells.Item(30, 3).Value = [double]"19.999999997853003"
$ws.Cells.Item(30, 4).Value = [double]"0.9684570266841195"
$ws.Cells.Item(30, 5).Value = [double]"0.24136208105394805"
$ws.Cells.Item(30, 6).Value = [double]"-0.040998626891935695"
$ws.Cells.Item(30, 7).Value = [double]"0.042622982887400775"
$ws.Cells.Item(31, 1).Value = [double]"259.61464508749606"
$ws.Cells.Item(31, 2).Value = [double]"248.32000131240653"
$ws.Cells.Item(31, 3).Value = [double]"19.999976078476294"
$ws.Cells.Item(31, 4).Value = [double]"0.9712028571171406"
$ws.Cells.Item(31, 5).Value = [double]"0.2361743106754203"
$ws.Cells.Item(31, 6).Value = [double]"-0.021989533480778294"
$ws.Cells.Item(31, 7).Value = [double]"0.015208081299252545"
$ws.Cells.Item(32, 1).Value = [double]"259.6060212463125"
$ws.Cells.Item(32, 2).Value = [double]"248.22206206968752"
$ws.Cells.Item(32, 3).Value = [double]"19.99997379700073"
$ws.Cells.Item(32, 4).Value = [double]"0.9712431062492859"
$ws.Cells.Item(32, 5).Value = [double]"0.23616598412877352"
$ws.Cells.Item(32, 6).Value = [double]"-0.02080418846565598"
$ws.Cells.Item(32, 7).Value = [double]"0.01444452144893002"
$ws.Cells.Item(33, 1).Value = [double]"259.59789432490425"
$ws.Cells.Item(33, 2).Value = [double]"248.124127855654"
$ws.Cells.Item(33, 3).Value = [double]"19.99997378720053"
$ws.Cells.Item(33, 4).Value = [double]"0.9712809929059415"
$ws.Cells.Item(33, 5).Value = [double]"0.23615877751817385"
$ws.Cells.Item(33, 6).Value = [double]"-0.019605890835138376"
$ws.Cells.Item(33, 7).Value = [double]"0.01370034916322505"
$ws.Cells.Item(34, 1).Value = [double]"259.5905504063567"
$ws.Cells.Item(34, 2).Value = [double]"248.0262203628549"
$ws.Cells.Item(34, 3).Value = [double]"19.999978298408372"
$ws.Cells.Item(34, 4).Value = [double]"0.9713604404898503"
$ws.Cells.Item(34, 5).Value = [double]"0.2361064226816632"
$ws.Cells.Item(34, 6).Value = [double]"-0.017717941598074223"
$ws.Cells.Item(34, 7).Value = [double]"0.011574708767048794"
